$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.768.34"
$ws.Range("E2").Value = "  +7.58%  "
$ws.Range("D3").Value = "3.546.00"
$ws.Range("E3").Value = "  +10.22%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "192.35"
$ws.Range("E5").Value = "  +10.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "555.41"
$ws.Range("E6").Value = "  +8.04%  "
$ws.Range("D7").Value = "3.539.56"
$ws.Range("E7").Value = "  +10.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.610"
$ws.Range("E8").Value = "  +3.52%  "
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.643"
$ws.Range("E10").Value = "  +7.86%  "
$ws.Range("E11").Value = "  +8.39%  "
$ws.Range("E12").Value = "  +16.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000274"
$ws.Range("E13").Value = "  +9.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.50"
$ws.Range("E14").Value = "  +7.25%  "
$ws.Range("D15").Value = "4.105.52"
$ws.Range("E15").Value = "  +10.26%  "
$ws.Range("D16").Value = "3.547.48"
$ws.Range("E16").Value = "  +10.56%  "
$ws.Range("D17").Value = "67.823.28"
$ws.Range("E17").Value = "  +8.01%  "
$ws.Range("E18").Value = "  +5.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.44"
$ws.Range("E19").Value = "  +7.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.96"
$ws.Range("E20").Value = "  +9.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  +5.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "408.44"
$ws.Range("E22").Value = "  +11.96%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.98"
$ws.Range("E23").Value = "  +7.53%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.29"
$ws.Range("E24").Value = "  +10.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.30"
$ws.Range("E25").Value = "  +6.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.46"
$ws.Range("E26").Value = "  +3.93%  "
$ws.Range("E27").Value = "  +14.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.16"
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "12.01"
$ws.Range("E29").Value = "  +7.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.84"
$ws.Range("E30").Value = "  +9.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "699.86"
$ws.Range("E31").Value = "  +7.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "30.64"
$ws.Range("E32").Value = "  +8.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.86"
$ws.Range("E33").Value = "  +10.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.84"
$ws.Range("E34").Value = "  +7.07%  "
$ws.Range("E35").Value = "  +8.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "60.74"
$ws.Range("E36").Value = "  +5.07%  "
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "39.32"
$ws.Range("E37").Value = "  +8.22%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0835"
$ws.Range("E38").Value = "  +20.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.399"
$ws.Range("E40").Value = "  +7.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.139"
$ws.Range("E41").Value = "  +14.13%  "
$ws.Range("E42").Value = "  +18.37%  "
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.01"
$ws.Range("E44").Value = "  +17.10%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "3.057.93"
$ws.Range("E45").Value = "  +6.71%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.70"
$ws.Range("E46").Value = "  +6.31%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0424"
$ws.Range("E47").Value = "  +9.45%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.29"
$ws.Range("E48").Value = "  +15.39%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.14"
$ws.Range("E49").Value = "  +20.98%  "
$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.73"
$ws.Range("E50").Value = "  +2.31%  "
$ws.Range("E51").Value = "  +7.43%  "
